# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) cells on the zh-cn and de-de
# sheets with newer handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 16:46:46"
$wsZh.Range("H2").Value = "2016-03-19 16:47:27"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 16:46:54"
$wsDe.Range("H2").Value = "2016-03-19 16:47:40"
